$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("customer_acc")
$ws2 = $wb.Worksheets.Item("customer_prod")
$ws3 = $wb.Worksheets.Item("valid_product_names")

# ---------------------------------------------------------------------------
# Sheet3 (valid_product_names): insert a header row, add Kenya/Zimbabwe lists
# ---------------------------------------------------------------------------
$ws3.Rows.Item(1).Insert()

$ws3.Range("A1").Value = "Kenya"
$ws3.Range("B1").Value = "Zimbabwe"

$ws3.Range("B2").Value = "Fixed Interest - Money Market"
$ws3.Range("B3").Value = "Equity - General"
$ws3.Range("B4").Value = "Real Estate - General"
$ws3.Range("B5").Value = "Flexi Funeral Plan"
$ws3.Range("B6").Value = "Funeral Plan"
$ws3.Range("B7").Value = "Term Plan"
$ws3.Range("B8").Value = "Life Plan"
$ws3.Range("B9").Value = "Savings Plan"

$ws3.Columns.Item(1).ColumnWidth = 37
$ws3.Columns.Item(2).ColumnWidth = 28.28515625

# ---------------------------------------------------------------------------
# workbook-level defined names for the country product lists
# ---------------------------------------------------------------------------
$wb.Names.Add("Kenya", "=valid_product_names!`$A`$2:`$A`$62")
$wb.Names.Add("Zimbabwe", "=valid_product_names!`$B`$2:`$B`$9")

# ---------------------------------------------------------------------------
# Sheet1 (customer_acc): new account numbers, drop the trailing blank row
# ---------------------------------------------------------------------------
$ws1.Rows.Item(20).Delete()
$ws1.Range("A4").Value = 100003
$ws1.Range("A5").Value = 101403

# ---------------------------------------------------------------------------
# Sheet2 (customer_prod): new country picker column + refreshed sample data
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Old Mutual Balanced Fund"
$ws2.Range("A3").Value = "CORPORATE"
$ws2.Range("A4").Value = "MOTOR PRIVATE"
$ws2.Range("A5").ClearContents()

# extend the used range down to row 19, keeping the same cell style as A16
$ws2.Range("A16").Copy($ws2.Range("A17:A19"))

# new country picker column, matching styles of the PRODUCT column
$ws2.Range("A1").Copy($ws2.Range("C1"))
$ws2.Range("C1").Value = "SELECT A COUNTRY"
$ws2.Range("A2").Copy($ws2.Range("C2"))
$ws2.Range("C2").Value = "Kenya"

$ws2.Columns.Item(2).ColumnWidth = 9.7109375
$ws2.Columns.Item(3).ColumnWidth = 18.140625

$ws2.Range("A2:A16").Validation.Delete()
$ws2.Range("A2:A19").Validation.Add(3, 1, 1, "=INDIRECT(`$C`$2)")
$ws2.Range("C2").Validation.Add(3, 1, 1, "=valid_product_names!`$A`$1:`$B`$1")

# ---------------------------------------------------------------------------
# Selections - sheet2 (customer_prod) must be selected last to stay the
# active tab, matching the workbook's saved state.
# ---------------------------------------------------------------------------
$ws1.Range("G18").Select()
$ws3.Range("B16").Select()
$ws2.Range("D4").Select()
